# Sprint Backlog.xlsx update:
#  - Add new "Actual" column (E)
#  - Add two new backlog rows (Favorite Seller + Radio Button filter tasks)
#  - Move/refresh the totals row and formulas accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

# ------------------------------------------------------------------
# 1. Insert two new rows before the current totals row (row 10),
#    pushing the totals row down to row 12 (and shifting the
#    existing merged cell A10:B10 -> A12:B12, and the existing
#    formatting on that row, automatically).
# ------------------------------------------------------------------
$ws.Range("A10:A11").EntireRow.Insert()

# ------------------------------------------------------------------
# 2. Fill in the "Actual" values for the existing task rows (2-9),
#    copying the formatting already used by the "Estimate" column (D)
#    so no stray new styles get introduced.
# ------------------------------------------------------------------
$ws.Range("D2:D9").Copy()
$ws.Range("E2:E9").PasteSpecial(-4122)

$actualValues = @(0.1, 0.1, 0.1, 0.1, 0.1, 1, 2, 1)
for ($i = 0; $i -lt $actualValues.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 5).Value2 = $actualValues[$i]
}

# ------------------------------------------------------------------
# 3. New backlog rows (10 & 11) - start from the previous task row's
#    formatting, then tint the whole new rows yellow.
# ------------------------------------------------------------------
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("A9:D9").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)

$ws.Range("A10").Value2 = "Implement Favorite Seller functionality"
$ws.Range("B10").Value2 = "Favorite Seller"
$ws.Range("C10").Value2 = "Brianna"
$ws.Range("E10").Value2 = 1

$ws.Range("A11").Value2 = "Implement Radion Button Functionality for filter"
$ws.Range("B11").Value2 = "Filter Discounts "
$ws.Range("C11").Value2 = "Destiny"
$ws.Range("E11").Value2 = 1.5

$ws.Range("A10:E11").Interior.Color = 65535

# ------------------------------------------------------------------
# 4. Header row - add "Actual" header in E1 (matching D1's base
#    formatting) and center B1:E1.
# ------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value2 = "Actual"

$ws.Range("B1:E1").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 5. Totals row (now row 12) - refresh formulas to cover new rows.
#    E12 reuses D12's existing formatting.
# ------------------------------------------------------------------
$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122)

$ws.Range("D12").Formula = "=SUM(D2:D10)"
$ws.Range("E12").Formula = "=SUM(E2:E11)"

# ------------------------------------------------------------------
# 6. Column width for new column E, and refresh the selection.
# ------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 8.6640625
$ws.Range("E14").Select()
